$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 7495
$ws.Range("I51").Value = 7495
$ws.Range("K51").Value = 7495
$ws.Range("M51").Value = -7011
# Row 55
$ws.Range("H55").Value = 2468.111
$ws.Range("I55").Value = 287.5
$ws.Range("J55").Value = 4212.6
$ws.Range("K55").Value = 287.5
$ws.Range("L55").Value = 4212.6
$ws.Range("M55").Value = -73.5
$ws.Range("N55").Value = -4640.6
# Row 62
$ws.Range("H62").Value = 5361.0415
$ws.Range("I62").Value = 3730.1052
$ws.Range("J62").Value = 11558.6
$ws.Range("K62").Value = 3730.1052
$ws.Range("L62").Value = 11558.6
$ws.Range("M62").Value = -3106.1052
$ws.Range("N62").Value = -12806.6
# Row 65
$ws.Range("H65").Value = 5361.0415
$ws.Range("I65").Value = 3730.1052
$ws.Range("J65").Value = 11558.6
$ws.Range("K65").Value = 18650.526
$ws.Range("L65").Value = 57793
$ws.Range("M65").Value = -15530.526
$ws.Range("N65").Value = -64033
# Row 93
$ws.Range("H93").Value = 54902
$ws.Range("J93").Value = 54902
$ws.Range("L93").Value = 54902
$ws.Range("N93").Value = -59894
# Row 103
$ws.Range("H103").Value = 2066.3333
$ws.Range("I103").Value = 2066.3333
$ws.Range("K103").Value = 6198.999899999999
$ws.Range("M103").Value = -5612.999899999999
# Row 112
$ws.Range("H112").Value = 27799.75
$ws.Range("I112").Value = 1199
$ws.Range("J112").Value = 36666.668
$ws.Range("K112").Value = 3597
$ws.Range("L112").Value = 110000.004
$ws.Range("N112").Value = -112216.004
$ws.Range("M112").Value = -2489
# Row 132
$ws.Range("H132").Value = 15153735
$ws.Range("I132").Value = 15627239
$ws.Range("K132").Value = 46881717
$ws.Range("M132").Value = -46879187
# Row 137
$ws.Range("H137").Value = 116829.56
$ws.Range("I137").Value = 227159.75
$ws.Range("K137").Value = 681479.25
$ws.Range("M137").Value = -678929.25
# Row 138
$ws.Range("H138").Value = 2735.6
$ws.Range("I138").Value = 751.48
$ws.Range("J138").Value = 3396.9734
$ws.Range("K138").Value = 2254.44
$ws.Range("L138").Value = 10190.9202
$ws.Range("M138").Value = 2885.56
$ws.Range("N138").Value = -20470.9202

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1535.08
$ws.Range("I32").Value = 1252.6063
$ws.Range("K32").Value = 1252.6063
$ws.Range("M32").Value = -965.6062999999999
# Row 61
$ws.Range("H61").Value = 1614.8948
$ws.Range("I61").Value = 1563.0769
$ws.Range("J61").Value = 1727.1666
$ws.Range("K61").Value = 1563.0769
$ws.Range("L61").Value = 1727.1666
$ws.Range("M61").Value = -1351.0769
$ws.Range("N61").Value = -2151.1666
# Row 132
$ws.Range("H132").Value = 1462.5491
$ws.Range("I132").Value = 1200.738
$ws.Range("K132").Value = 3602.214
$ws.Range("M132").Value = -1072.214
# Row 136
$ws.Range("H136").Value = 1614.8948
$ws.Range("I136").Value = 1563.0769
$ws.Range("J136").Value = 1727.1666
$ws.Range("K136").Value = 4689.2307
$ws.Range("L136").Value = 5181.4998
$ws.Range("M136").Value = -2139.2307
$ws.Range("N136").Value = -10281.4998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2904.4783
$ws.Range("I134").Value = 1191.2188
$ws.Range("K134").Value = 3573.6564
$ws.Range("M134").Value = -1038.6564

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 396.41177
$ws.Range("J7").Value = 739.4
$ws.Range("L7").Value = 739.4
$ws.Range("N7").Value = -965.4
# Row 31
$ws.Range("H31").Value = 31483
$ws.Range("I31").Value = 1916
$ws.Range("K31").Value = 1916
$ws.Range("M31").Value = -1621
# Row 34
$ws.Range("H34").Value = 31483
$ws.Range("I34").Value = 1916
$ws.Range("K34").Value = 1916
$ws.Range("M34").Value = -1714
# Row 62
$ws.Range("H62").Value = 6966.5
$ws.Range("I62").Value = 1998
$ws.Range("J62").Value = 11935
$ws.Range("K62").Value = 1998
$ws.Range("L62").Value = 11935
$ws.Range("M62").Value = -1374
$ws.Range("N62").Value = -13183
# Row 65
$ws.Range("H65").Value = 6966.5
$ws.Range("I65").Value = 1998
$ws.Range("J65").Value = 11935
$ws.Range("K65").Value = 9990
$ws.Range("L65").Value = 59675
$ws.Range("M65").Value = -6870
$ws.Range("N65").Value = -65915
# Row 132
$ws.Range("H132").Value = 19969.66
$ws.Range("I132").Value = 1968.75
$ws.Range("K132").Value = 5906.25
$ws.Range("M132").Value = -3376.25
# Row 134
$ws.Range("H134").Value = 2811.1206
$ws.Range("I134").Value = 2409.0205
$ws.Range("J134").Value = 5000.3335
$ws.Range("K134").Value = 7227.0615
$ws.Range("L134").Value = 15001.0005
$ws.Range("M134").Value = -4692.0615
$ws.Range("N134").Value = -20071.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 47529.9
$ws.Range("J37").Value = 47529.9
$ws.Range("L37").Value = 142589.7
$ws.Range("N37").Value = -142813.7
# Row 56
$ws.Range("H56").Value = 20839654
$ws.Range("I56").Value = 20839654
$ws.Range("K56").Value = 20839654
$ws.Range("M56").Value = -20839124
# Row 122
$ws.Range("H122").Value = 1194.1538
$ws.Range("J122").Value = 1325.8334
$ws.Range("L122").Value = 11932.5006
$ws.Range("N122").Value = -16832.5006
# Row 123
$ws.Range("H123").Value = 4032.25
$ws.Range("I123").Value = 4032.25
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 12096.75
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -9646.75
$ws.Range("N123").Value = ""
# Row 124
$ws.Range("H124").Value = 1997.5
$ws.Range("I124").Value = 1995
$ws.Range("K124").Value = 5985
$ws.Range("M124").Value = -1075
# Row 125
$ws.Range("H125").Value = 2449.5
$ws.Range("I125").Value = 2449.5
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 7348.5
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -2428.5
$ws.Range("N125").Value = ""
# Row 139
$ws.Range("H139").Value = 1477.8
$ws.Range("J139").Value = 2733
$ws.Range("L139").Value = 8199
$ws.Range("N139").Value = -18479

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3575.9333
$ws.Range("I132").Value = 3681.32
$ws.Range("K132").Value = 11043.96
$ws.Range("M132").Value = -8513.960000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8219
$ws.Range("I7").Value = 6765.778
$ws.Range("J7").Value = 10398.833
$ws.Range("K7").Value = 6765.778
$ws.Range("L7").Value = 10398.833
$ws.Range("M7").Value = -6653.778
$ws.Range("N7").Value = -10622.833
# Row 22
$ws.Range("H22").Value = 297212.66
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("N22").Value = -2590
# Row 27
$ws.Range("H27").Value = 297212.66
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 2000
$ws.Range("N27").Value = -2214
# Row 68
$ws.Range("H68").Value = 5054.375
$ws.Range("J68").Value = 5499.25
$ws.Range("L68").Value = 5499.25
$ws.Range("N68").Value = -6997.25
# Row 71
$ws.Range("H71").Value = 5054.375
$ws.Range("J71").Value = 5499.25
$ws.Range("L71").Value = 27496.25
$ws.Range("N71").Value = -34984.25
# Row 126
$ws.Range("H126").Value = 8219
$ws.Range("I126").Value = 6765.778
$ws.Range("J126").Value = 10398.833
$ws.Range("K126").Value = 20297.334
$ws.Range("L126").Value = 31196.499
$ws.Range("M126").Value = -17827.334
$ws.Range("N126").Value = -36136.499
# Row 132
$ws.Range("H132").Value = 7029.971
$ws.Range("I132").Value = 7896.52
$ws.Range("J132").Value = 4863.6
$ws.Range("K132").Value = 23689.56
$ws.Range("L132").Value = 14590.8
$ws.Range("M132").Value = -21159.56
$ws.Range("N132").Value = -19650.8
# Row 136
$ws.Range("H136").Value = 61522.61
$ws.Range("I136").Value = 170250.5
$ws.Range("J136").Value = 7158.6665
$ws.Range("K136").Value = 510751.5
$ws.Range("L136").Value = 21475.9995
$ws.Range("M136").Value = -508201.5
$ws.Range("N136").Value = -26575.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 2950
$ws.Range("I2").Value = 2900
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 2900
$ws.Range("L2").Value = 3000
$ws.Range("N2").Value = -3224
$ws.Range("M2").Value = -2788
# Row 62
$ws.Range("H62").Value = 8820.120000000001
$ws.Range("I62").Value = 4100
$ws.Range("J62").Value = 9230.565000000001
$ws.Range("K62").Value = 4100
$ws.Range("L62").Value = 9230.565000000001
$ws.Range("M62").Value = -3476
$ws.Range("N62").Value = -10478.565
# Row 65
$ws.Range("H65").Value = 8820.120000000001
$ws.Range("I65").Value = 4100
$ws.Range("J65").Value = 9230.565000000001
$ws.Range("K65").Value = 20500
$ws.Range("L65").Value = 46152.825
$ws.Range("M65").Value = -17380
$ws.Range("N65").Value = -52392.825
# Row 81
$ws.Range("H81").Value = 27779378
$ws.Range("I81").Value = 27779378
$ws.Range("K81").Value = 55558756
$ws.Range("M81").Value = -55557695
# Row 84
$ws.Range("H84").Value = 27779378
$ws.Range("I84").Value = 27779378
$ws.Range("K84").Value = 277793780
$ws.Range("M84").Value = -277788476
# Row 107
$ws.Range("H107").Value = 66667490
$ws.Range("J107").Value = 987.5
$ws.Range("L107").Value = 2962.5
$ws.Range("N107").Value = -6802.5
# Row 132
$ws.Range("H132").Value = 50051320
$ws.Range("I132").Value = 62507964
$ws.Range("K132").Value = 187523892
$ws.Range("M132").Value = -187521362
# Row 136
$ws.Range("H136").Value = 2838.9524
$ws.Range("I136").Value = 2050.2856
$ws.Range("J136").Value = 4416.2856
$ws.Range("K136").Value = 6150.8568
$ws.Range("L136").Value = 13248.8568
$ws.Range("M136").Value = -3600.8568
$ws.Range("N136").Value = -18348.8568
